$wb = $excel.ActiveWorkbook

# --- Sheet "Salesforce": insert new "LeadName" column (with "Anil" as the lead
#     name value) right before the existing "QuoteType" column, and fix the
#     typo "YearABc" -> "YearABC" while we're in there. ---
$ws = $wb.Worksheets.Item("Salesforce")
$ws.Activate()

$ws.Columns("O").Insert()
$ws.Range("O1").Value = "LeadName"
$ws.Range("O2").Value = "Anil"

# typo fix picked up in the same commit
$ws.Range("L2").Value = "YearABC"

# The column insert (in this sandbox) does not carry existing hyperlink
# anchors along with the shifted cells, so rebuild all of this sheet's
# hyperlinks at their correct, post-insert locations.
$ws.Range("A1:Z10").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "https://testingxperts-5d-dev-ed.develop.my.salesforce.com", "", "", "https://testingxperts-5d-dev-ed.develop.my.salesforce.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:xperttesting3@gmail.com", "", "", "xperttesting3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Tiger@2024", "", "", "Tiger@2024")
$ws.Hyperlinks.Add($ws.Range("Q2"), "mailto:demo@yahooo.com", "", "", "demo@yahooo.com")

# Selection ends up sitting on the newly typed header cell.
$ws.Range("O1").Select()
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1

# --- Switch the active sheet to "SalesforceLogin" (last sheet), which is
#     where the user ended up / saved from. ---
$ws8 = $wb.Worksheets.Item("SalesforceLogin")
$ws8.Activate()
